$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update part row 11: SM2V05 (0.31" travel coupler) -> SM2V10 (0.81" travel coupler)
$ws.Cells.Item(11, 1).Value2 = "SM2V10"
$ws.Cells.Item(11, 3).Value2 = "Ø2`" Adjustable Lens Tube, 0.81`" Travel"
$ws.Cells.Item(11, 5).Value2 = 54
$ws.Cells.Item(11, 7).Value2 = "SM2 coupler to mount and rotate the turret relative to filter wheel "

# Move the active selection to C12, matching the saved view state
$ws.Range("C12").Select()

# Reposition/resize the workbook window to match the saved view state
$win = $wb.Windows.Item(1)
$win.Left = 5376
$win.Top = 5160
$win.Width = 34560
$win.Height = 18684
